$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J header
$ws.Cells.Item(1, 10).Value = "portrait"

# Rename NAPOLEON 1er -> Napoleon
$ws.Cells.Item(27, 2).Value = "Napoleon"

# Portrait file names for the "generaux" rows (French army)
$ws.Cells.Item(27, 10).Value = "Napoleon.webp"
$ws.Cells.Item(28, 10).Value = "Bernadotte.jpg"
$ws.Cells.Item(29, 10).Value = "Davout.jpg"
$ws.Cells.Item(30, 10).Value = "Soult.jpg"
$ws.Cells.Item(31, 10).Value = "Lannes.jpg"
$ws.Cells.Item(32, 10).Value = "Murat.jpg"

# Portrait file names for the "generaux" rows (Austro-Russian army)
$ws.Cells.Item(63, 10).Value = "AlexandreIer.jpg"
$ws.Cells.Item(64, 10).Value = "Buxhoevden.jpg"
$ws.Cells.Item(65, 10).Value = "Koutouzof.jpg"
$ws.Cells.Item(66, 10).Value = "Bagration.jpg"
$ws.Cells.Item(67, 10).Value = "Liechtenstein.jpg"
$ws.Cells.Item(68, 10).Value = "RussiaFlag.jpg"

# Widen new column J and set the active cell roughly where the author left off
$ws.Columns.Item(10).ColumnWidth = 13
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$ws.Application.Goto($ws.Range("A43"), $true) | Out-Null
$ws.Range("L65").Select() | Out-Null
